# Generate Report for Handoff
# The localization-status report was regenerated a bit later than before:
# several "latest handoff/handback" timestamps that were stamped at
# 18:21:xx now collapse to the newer run's timestamps (18:21:58 / 18:22:02
# depending on locale), since this run happened slightly later than the
# previous one.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column D ("Latest Handoff Date") ---
$ws = $wb.Worksheets.Item("Overview")
$newOverviewDate = "2016-22-17 18:22:02"
foreach ($r in 7,10,11,12,13,14,15,16) {
    $ws.Cells.Item($r, 4).Value = $newOverviewDate
}

# --- zh-cn sheet: column E ("Latest Handoff Datetime") ---
$ws = $wb.Worksheets.Item("zh-cn")
$newZhDate = "2016-03-17 18:21:58"
foreach ($r in 7,10,11,12,13,14,15,16) {
    $ws.Cells.Item($r, 5).Value = $newZhDate
}

# --- de-de sheet: column E ("Latest Handoff Datetime") ---
$ws = $wb.Worksheets.Item("de-de")
$newDeDate = "2016-03-17 18:22:02"
foreach ($r in 7,10,11,12,13,14,15,16) {
    $ws.Cells.Item($r, 5).Value = $newDeDate
}
